# AzureAD master application and setDataFromDataExcel
#
# Sheet "AZUREAD" (sheet1): extend the "Tabla1" table with 3 new columns
# (URL, CLIENT ID, CLIENT SECRET) and populate the 3 data rows.
# Sheet "AD" (sheet2): populate the single data row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "AZUREAD"
# ---------------------------------------------------------------------
$wsAzure = $wb.Worksheets.Item("AZUREAD")
$tbl = $wsAzure.ListObjects.Item(1)

# Grow the table from 3 to 6 columns (A:C -> A:F); new columns get
# placeholder names that are overwritten below via the header cells.
$tbl.ListColumns.Add() | Out-Null
$tbl.ListColumns.Add() | Out-Null
$tbl.ListColumns.Add() | Out-Null

# --- Header row (row 1) ------------------------------------------------
$wsAzure.Range("A1").Value = "TYPE"
$wsAzure.Range("B1").Value = "SCOPE"
$wsAzure.Range("C1").Value = "APP NAME"
$wsAzure.Range("D1").Value = "URL"
$wsAzure.Range("E1").Value = "CLIENT ID"
$wsAzure.Range("F1").Value = "CLIENT SECRET"

# Re-apply the header style (centered, like A1:C1) to the new header
# cells D1:F1 by copying the format from an existing header cell.
$wsAzure.Range("A1").Copy() | Out-Null
$wsAzure.Range("D1:F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Data rows -----------------------------------------------------------
$wsAzure.Range("A2").Value = "AzureAD"
$wsAzure.Range("B2").Value = "SGTO"
$wsAzure.Range("C2").Value = "Celonis"
$wsAzure.Range("D2").Value = "-"
$wsAzure.Range("E2").Value = "-"
$wsAzure.Range("F2").Value = "-"

$wsAzure.Range("A3").Value = "AzureAD"
$wsAzure.Range("B3").Value = "SGTO"
$wsAzure.Range("C3").Value = "O365"
$wsAzure.Range("D3").Value = "-"
$wsAzure.Range("E3").Value = "-"
$wsAzure.Range("F3").Value = "-"

$wsAzure.Range("A4").Value = "AzureAD"
$wsAzure.Range("B4").Value = "SGTO"
$wsAzure.Range("C4").Value = "Irisrusk"
$wsAzure.Range("D4").Value = "URL"
$wsAzure.Range("E4").Value = "CLIENTID"
$wsAzure.Range("F4").Value = "CLIENTSECRET"

# Column widths for the 3 new columns (D:F)
$wsAzure.Range("D1:F1").ColumnWidth = 29.83

# Selection ends on F3 after the edits
$wsAzure.Range("F3").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "AD"
# ---------------------------------------------------------------------
$wsAD = $wb.Worksheets.Item("AD")

$wsAD.Range("A1").Value = "TYPE"
$wsAD.Range("B1").Value = "SCOPE"
$wsAD.Range("C1").Value = "APP NAME"

$wsAD.Range("A2").Value = "AD"
$wsAD.Range("B2").Value = "SGTO"
$wsAD.Range("C2").Value = "VPNs"

$wsAD.Range("B2").Select() | Out-Null

# Leave the workbook with "AZUREAD" as the active/visible tab and the
# selection on F3, matching the final state after the edits.
$wsAzure.Activate() | Out-Null
$wsAzure.Range("F3").Select() | Out-Null
